$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-05 Sunday" "2025-10-06 Monday"

Replace-Text "94×22=" "34×25="
Replace-Text "63×91=" "82×53="
Replace-Text "68×52=" "17×38="
Replace-Text "65×87=" "80×22="
Replace-Text "41×30=" "27×30="

Replace-Text "97×62=" "61×88="
Replace-Text "37×96=" "50×73="
Replace-Text "63×52=" "70×11="
Replace-Text "32×39=" "67×74="
Replace-Text "70×52=" "58×70="

Replace-Text "76×14=" "88×49="
Replace-Text "87×42=" "73×65="
Replace-Text "57×13=" "26×56="
Replace-Text "24×47=" "92×83="
Replace-Text "18×64=" "21×93="

Replace-Text "70×75=" "77×17="
Replace-Text "57×34=" "31×80="
Replace-Text "86×51=" "75×75="
Replace-Text "53×42=" "32×50="
Replace-Text "83×68=" "74×80="

Replace-Text "65×86=" "84×73="
Replace-Text "35×78=" "79×68="
Replace-Text "54×95=" "21×58="
Replace-Text "50×87=" "53×90="
Replace-Text "41×99=" "30×34="
